$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 1288.3334
$ws.Range("J4").Value = 500
$ws.Range("L4").Value = 500
$ws.Range("N4").Value = -728
$ws.Range("H63").Value = 50102.6
$ws.Range("I63").Value = 50000
$ws.Range("J63").Value = 50128.25
$ws.Range("K63").Value = 50000
$ws.Range("L63").Value = 50128.25
$ws.Range("M63").Value = -49376
$ws.Range("N63").Value = -51376.25
$ws.Range("H66").Value = 50102.6
$ws.Range("I66").Value = 50000
$ws.Range("J66").Value = 50128.25
$ws.Range("K66").Value = 150000
$ws.Range("L66").Value = 150384.75
$ws.Range("M66").Value = -146880
$ws.Range("N66").Value = -156624.75
$ws.Range("H74").Value = 3356.08
$ws.Range("I74").Value = 3200.1667
$ws.Range("J74").Value = 3500
$ws.Range("K74").Value = 3200.1667
$ws.Range("L74").Value = 3500
$ws.Range("M74").Value = -2264.1667
$ws.Range("N74").Value = -5372
$ws.Range("H76").Value = 6502.923
$ws.Range("I76").Value = 4626.5713
$ws.Range("J76").Value = 8692
$ws.Range("K76").Value = 4626.5713
$ws.Range("L76").Value = 8692
$ws.Range("M76").Value = -4311.5713
$ws.Range("N76").Value = -9322
$ws.Range("H77").Value = 3356.08
$ws.Range("I77").Value = 3200.1667
$ws.Range("J77").Value = 3500
$ws.Range("K77").Value = 16000.8335
$ws.Range("L77").Value = 17500
$ws.Range("M77").Value = -11320.8335
$ws.Range("N77").Value = -26860
$ws.Range("H79").Value = 6502.923
$ws.Range("I79").Value = 4626.5713
$ws.Range("J79").Value = 8692
$ws.Range("K79").Value = 4626.5713
$ws.Range("L79").Value = 8692
$ws.Range("M79").Value = -3534.5713
$ws.Range("N79").Value = -10876
$ws.Range("H138").Value = 7577624.5
$ws.Range("I138").Value = 1700.5416
$ws.Range("J138").Value = 27780088
$ws.Range("K138").Value = 5101.6248
$ws.Range("L138").Value = 83340264
$ws.Range("M138").Value = 38.3752000000004
$ws.Range("N138").Value = -83350544

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 144.5
$ws.Range("I4").Value = 92.666664
$ws.Range("K4").Value = 92.666664
$ws.Range("M4").Value = 23.333336
$ws.Range("H32").Value = 9959.082
$ws.Range("I32").Value = 10378.953
$ws.Range("K32").Value = 10378.953
$ws.Range("M32").Value = -10091.953
$ws.Range("H63").Value = 166667970
$ws.Range("J63").Value = 1948
$ws.Range("L63").Value = 1948
$ws.Range("N63").Value = -3320
$ws.Range("H66").Value = 166667970
$ws.Range("J66").Value = 1948
$ws.Range("L66").Value = 9740
$ws.Range("N66").Value = -16604

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 444.375
$ws.Range("I22").Value = 445.83334
$ws.Range("J22").Value = 440
$ws.Range("K22").Value = 445.83334
$ws.Range("L22").Value = 440
$ws.Range("M22").Value = -272.83334
$ws.Range("N22").Value = -786
$ws.Range("H105").Value = 4579.967
$ws.Range("I105").Value = 3333.3333
$ws.Range("J105").Value = 4718.4814
$ws.Range("K105").Value = 3333.3333
$ws.Range("L105").Value = 4718.4814
$ws.Range("M105").Value = -1586.3333
$ws.Range("N105").Value = -8212.481400000001
$ws.Range("H134").Value = 3958.1936
$ws.Range("I134").Value = 3068.1707
$ws.Range("J134").Value = 5695.857
$ws.Range("K134").Value = 9204.5121
$ws.Range("L134").Value = 17087.571
$ws.Range("M134").Value = -6669.5121
$ws.Range("N134").Value = -22157.571

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 86
$ws.Range("I7").Value = 87.28570999999999
$ws.Range("J7").Value = 83.75
$ws.Range("K7").Value = 87.28570999999999
$ws.Range("L7").Value = 83.75
$ws.Range("M7").Value = 25.71429000000001
$ws.Range("N7").Value = -309.75
$ws.Range("H31").Value = 9527419
$ws.Range("I31").Value = 1990.65
$ws.Range("J31").Value = 22227990
$ws.Range("K31").Value = 1990.65
$ws.Range("L31").Value = 22227990
$ws.Range("M31").Value = -1695.65
$ws.Range("N31").Value = -22228580
$ws.Range("H34").Value = 9527419
$ws.Range("I34").Value = 1990.65
$ws.Range("J34").Value = 22227990
$ws.Range("K34").Value = 1990.65
$ws.Range("L34").Value = 22227990
$ws.Range("M34").Value = -1788.65
$ws.Range("N34").Value = -22228394
$ws.Range("H94").Value = 3417.6191
$ws.Range("I94").Value = 1147.8572
$ws.Range("K94").Value = 1147.8572
$ws.Range("M94").Value = -696.8571999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 18.818182
$ws.Range("I12").Value = 14
$ws.Range("K12").Value = 42
$ws.Range("M12").Value = 131
$ws.Range("H18").Value = 83334104
$ws.Range("I18").Value = 100000776
$ws.Range("J18").Value = 721
$ws.Range("K18").Value = 300002328
$ws.Range("L18").Value = 2163
$ws.Range("M18").Value = -300002159
$ws.Range("N18").Value = -2501
$ws.Range("H20").Value = 2775.7878
$ws.Range("I20").Value = 1800.1666
$ws.Range("J20").Value = 2992.5925
$ws.Range("K20").Value = 5400.4998
$ws.Range("L20").Value = 8977.7775
$ws.Range("M20").Value = -5173.4998
$ws.Range("N20").Value = -9431.7775
$ws.Range("H86").Value = 1161.375
$ws.Range("I86").Value = 939.1667
$ws.Range("J86").Value = 1383.5834
$ws.Range("K86").Value = 2817.5001
$ws.Range("L86").Value = 4150.7502
$ws.Range("M86").Value = -1631.5001
$ws.Range("N86").Value = -6522.7502
$ws.Range("H89").Value = 1161.375
$ws.Range("I89").Value = 939.1667
$ws.Range("J89").Value = 1383.5834
$ws.Range("K89").Value = 8452.5003
$ws.Range("L89").Value = 12452.2506
$ws.Range("M89").Value = -2524.5003
$ws.Range("N89").Value = -24308.2506

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 15744.02
$ws.Range("I70").Value = 25456.465
$ws.Range("J70").Value = 4412.8335
$ws.Range("K70").Value = 25456.465
$ws.Range("L70").Value = 4412.8335
$ws.Range("M70").Value = -25186.465
$ws.Range("N70").Value = -4952.8335
$ws.Range("H73").Value = 15744.02
$ws.Range("I73").Value = 25456.465
$ws.Range("J73").Value = 4412.8335
$ws.Range("K73").Value = 25456.465
$ws.Range("L73").Value = 4412.8335
$ws.Range("M73").Value = -24520.465
$ws.Range("N73").Value = -6284.8335
$ws.Range("H80").Value = 15379484
$ws.Range("I80").Value = 22224650
$ws.Range("J80").Value = 3970875.2
$ws.Range("K80").Value = 22224650
$ws.Range("L80").Value = 3970875.2
$ws.Range("M80").Value = -22223652
$ws.Range("N80").Value = -3972871.2
$ws.Range("H83").Value = 15379484
$ws.Range("I83").Value = 22224650
$ws.Range("J83").Value = 3970875.2
$ws.Range("K83").Value = 111123250
$ws.Range("L83").Value = 19854376
$ws.Range("M83").Value = -111118258
$ws.Range("N83").Value = -19864360
$ws.Range("H113").Value = 1685.0714
$ws.Range("I113").Value = 1762.8182
$ws.Range("J113").Value = 1400
$ws.Range("K113").Value = 1762.8182
$ws.Range("L113").Value = 1400
$ws.Range("M113").Value = 407.1818000000001
$ws.Range("N113").Value = -5740

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2325.3
$ws.Range("I82").Value = 2410
$ws.Range("J82").Value = 2304.125
$ws.Range("K82").Value = 2410
$ws.Range("L82").Value = 2304.125
$ws.Range("M82").Value = -2049
$ws.Range("N82").Value = -3026.125
$ws.Range("H85").Value = 2325.3
$ws.Range("I85").Value = 2410
$ws.Range("J85").Value = 2304.125
$ws.Range("K85").Value = 2410
$ws.Range("L85").Value = 2304.125
$ws.Range("M85").Value = -1162
$ws.Range("N85").Value = -4800.125
